$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.463.08'
$ws.Range("E2").Value = '  +0.52%  '
$ws.Range("D3").Value = '3.096.51'
$ws.Range("E3").Value = '  -0.68%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.23'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.45%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '3.088.55'
$ws.Range("E8").Value = '  -0.65%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.527'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.37%  '
$ws.Range("E10").Value = '  +6.95%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.63'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.06%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.455'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.52%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000245'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.21'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.73%  '
$ws.Range("E15").Value = '  -1.19%  '
$ws.Range("D16").Value = '3.612.74'
$ws.Range("E16").Value = '  -0.68%  '
$ws.Range("D17").Value = '63.393.79'
$ws.Range("E17").Value = '  +0.51%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.08'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = '3.093.44'
$ws.Range("E19").Value = '  -0.86%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '461.28'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.15%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.21'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.721'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.70%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.42'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.38%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.22'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.95%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.86'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.25%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.13'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.59%  '
$ws.Range("E27").Value = '  +0.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.98'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +8.78%  '
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.66'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.57%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.19'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.73%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.88'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.59%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.110'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.22%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.60'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.42%  '
$ws.Range("D35").Value = '0.0₃0845'
$ws.Range("E35").Value = '  -2.56%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.41'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.35%  '
$ws.Range("E37").Value = '  -1.49%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.29'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.98%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.98'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.14%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '50.24'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.53%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '432.52'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.66'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.57%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0366'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.87%  '
$ws.Range("D44").Value = '2.874.14'
$ws.Range("E44").Value = '  -1.84%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.269'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.43%  '
$ws.Range("E46").Value = '  -3.32%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '35.45'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.10%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '123.90'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.21%  '
$ws.Range("E50").Value = '  -1.17%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.98'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.11%  '
